$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 4: "sarvn" user (copy of the "saravana" row with new values) ---
$ws.Range("A4").Value = "sarvn"
$ws.Range("B4").Value = "r"
$ws.Range("C4").Value = "sarvn.r"
$ws.Range("D4").Value = "OU=2D,OU=Spectrepost Users,DC=saravana,DC=com"
$ws.Range("G4").Value = "Production"

# E4 reuses the same password hyperlink as E3
$ws.Range("E4").Value = "D3skt0p@123"
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:D3skt0p@123", "", "", "D3skt0p@123") | Out-Null
$ws.Range("E4").Style = $ws.Range("E3").Style

# F4 gets its own new e-mail hyperlink
$ws.Range("F4").Value = "sarvn.r@saravana.com"
$ws.Hyperlinks.Add($ws.Range("F4"), "mailto:sarvn.r@saravana.com", "", "", "sarvn.r@saravana.com") | Out-Null

# --- Cosmetic / view-state tweaks ---
# move the active selection to the newly added cell, like the original author left it
$ws.Range("G4").Select() | Out-Null
